$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price/Volume columns as Text first so numeric-looking strings
# (e.g. "1.002") are preserved verbatim instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.061.59'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.835.01'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '244.75'
$ws.Range('E5').Value = '  +1.80%  '
$ws.Range('D6').Value = '0.6341'
$ws.Range('E6').Value = '  +2.46%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.07559'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').Value = '0.2948'
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('D10').Value = '22.92'
$ws.Range('E10').Value = '  +1.47%  '
$ws.Range('D11').Value = '0.07743'
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('D12').Value = '1.834.75'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '5.011'
$ws.Range('E13').Value = '  +1.36%  '
$ws.Range('D14').Value = '0.6725'
$ws.Range('E14').Value = '  +1.78%  '
$ws.Range('D15').Value = '83.41'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').Value = '0.000009608'
$ws.Range('E16').Value = '  +6.07%  '
$ws.Range('D17').Value = '6.110'
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').Value = '29.097.02'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = '12.61'
$ws.Range('E19').Value = '  +2.56%  '
$ws.Range('D20').Value = '227.13'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '7.209'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '160.60'
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('D25').Value = '0.1400'
$ws.Range('E25').Value = '  +3.74%  '
$ws.Range('D26').Value = '8.555'
$ws.Range('E26').Value = '  +2.09%  '
$ws.Range('D27').Value = '17.96'
$ws.Range('E27').Value = '  +0.98%  '
$ws.Range('D28').Value = '1.501'
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').Value = '4.126'
$ws.Range('E29').Value = '  +2.27%  '
$ws.Range('D30').Value = '4.078'
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').Value = '1.204'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').Value = '0.05390'
$ws.Range('E32').Value = '  +3.67%  '
$ws.Range('D33').Value = '1.869'
$ws.Range('E33').Value = '  +2.40%  '
$ws.Range('D34').Value = '0.7470'
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('D35').Value = '1.143'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').Value = '2.658'
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('D37').Value = '1.243.79'
$ws.Range('E37').Value = '  -2.18%  '
$ws.Range('D38').Value = '2.762'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('D39').Value = '0.01790'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').Value = '6.631'
$ws.Range('E40').Value = '  +5.04%  '
$ws.Range('D41').Value = '0.9073'
$ws.Range('E41').Value = '  +1.77%  '
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '102.17'
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('D44').Value = '1.986.98'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').Value = '0.00000000125'
$ws.Range('E45').Value = '  +5.37%  '
$ws.Range('D46').Value = '65.08'
$ws.Range('E46').Value = '  +2.93%  '
$ws.Range('D47').Value = '0.5116'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').Value = '0.4096'
$ws.Range('E48').Value = '  +3.77%  '
$ws.Range('D49').Value = '9.098'
$ws.Range('E49').Value = '  +2.59%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '1.656'
$ws.Range('E50').Value = '  -1.41%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = '6.781'
$ws.Range('E51').Value = '  +2.23%  '
